$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 (Product 50)
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 'Ceramic Pot - 8 inch'
$ws.Range("C51").Value = 'Durable ceramic pot perfect for indoor plants.'
$ws.Range("D51").Value = 50
$ws.Range("E51").Value = 15.99
$ws.Range("F51").Value = 10
$ws.Range("G51").Value = $true
$ws.Range("H51").Value = '[''All'', ''Pots'']'

# Row 52 (Product 51)
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 'Mini Pot Planters - 10 pcs'
$ws.Range("C52").Value = 'Set of 10 colorful mini planters for small plants.'
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = 12.5
$ws.Range("F52").Value = 8
$ws.Range("G52").Value = $true
$ws.Range("H52").Value = '[''All'', ''Pots'']'

# Row 53 (Product 52)
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 'Terracotta Pot - 12 inch'
$ws.Range("C53").Value = 'Classic terracotta pot for garden and patio use.'
$ws.Range("D53").Value = 40
$ws.Range("E53").Value = 18.75
$ws.Range("F53").Value = 12
$ws.Range("G53").Value = $true
$ws.Range("H53").Value = '[''All'', ''Pots'']'

# Row 54 (Product 53)
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 'Hanging Pot Set - 3 pcs'
$ws.Range("C54").Value = 'Hanging pot set with metal chains for balconies.'
$ws.Range("D54").Value = 30
$ws.Range("E54").Value = 22
$ws.Range("F54").Value = 15
$ws.Range("G54").Value = $true
$ws.Range("H54").Value = '[''All'', ''Pots'']'

# Row 55 (Product 54)
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 'Glass Pot for Succulents'
$ws.Range("C55").Value = 'Transparent glass pot ideal for succulent display.'
$ws.Range("D55").Value = 60
$ws.Range("E55").Value = 14.99
$ws.Range("F55").Value = 9.5
$ws.Range("G55").Value = $true
$ws.Range("H55").Value = '[''All'', ''Pots'']'

# Row 56 (Product 55)
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 'Large Outdoor Pot - 20 inch'
$ws.Range("C56").Value = 'Large outdoor pot designed for trees and shrubs.'
$ws.Range("D56").Value = 20
$ws.Range("E56").Value = 35
$ws.Range("F56").Value = 25
$ws.Range("G56").Value = $true
$ws.Range("H56").Value = '[''All'', ''Pots'']'

# Row 57 (Product 56)
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 'Decorative Pot with Stand'
$ws.Range("C57").Value = 'Decorative pot with wooden stand for living rooms.'
$ws.Range("D57").Value = 25
$ws.Range("E57").Value = 29.99
$ws.Range("F57").Value = 20
$ws.Range("G57").Value = $true
$ws.Range("H57").Value = '[''All'', ''Pots'']'

# Row 58 (Product 57)
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 'Plastic Pot - 5 inch'
$ws.Range("C58").Value = 'Sturdy plastic pot suitable for various plants.'
$ws.Range("D58").Value = 70
$ws.Range("E58").Value = 10
$ws.Range("F58").Value = 6.5
$ws.Range("G58").Value = $true
$ws.Range("H58").Value = '[''All'', ''Pots'']'

# Row 59 (Product 58)
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 'Bamboo Pot - Eco Friendly'
$ws.Range("C59").Value = 'Eco-friendly bamboo pot for sustainable gardening.'
$ws.Range("D59").Value = 45
$ws.Range("E59").Value = 16.5
$ws.Range("F59").Value = 11
$ws.Range("G59").Value = $true
$ws.Range("H59").Value = '[''All'', ''Pots'']'

# Row 60 (Product 59)
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 'Self-Watering Pot - 6 inch'
$ws.Range("C60").Value = 'Self-watering pot to keep plants hydrated longer.'
$ws.Range("D60").Value = 55
$ws.Range("E60").Value = 19.99
$ws.Range("F60").Value = 13
$ws.Range("G60").Value = $true
$ws.Range("H60").Value = '[''All'', ''Pots'']'
